# LBCB 3 Calibration settings.
# Adds the new "Offsets" summary table (LVDT / Servo Error columns) on Sheet3
# and makes Sheet3 the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

# Column widths for the new table
$ws3.Columns.Item(1).ColumnWidth = 12.02
$ws3.Columns.Item(2).ColumnWidth = 19.88
$ws3.Columns.Item(3).ColumnWidth = 17.88

# Header row
$ws3.Range("A1").Value = "Offsets"
$ws3.Range("B1").Value = "LVDT"
$ws3.Range("C1").Value = "Servo Error"

# Data rows
$ws3.Range("A2").Value = "X1"
$ws3.Range("B2").Value = 5.0604899999999997
$ws3.Range("C2").Value = -0.0142364

$ws3.Range("A3").Value = "X2"
$ws3.Range("B3").Value = 5.0463800000000001
$ws3.Range("C3").Value = 0.0119929

$ws3.Range("A4").Value = "Y1"
$ws3.Range("B4").Value = 4.9553200000000004
$ws3.Range("C4").Value = 0.0116778

$ws3.Range("A5").Value = "Z1"
$ws3.Range("B5").Value = 4.8736199999999998
$ws3.Range("C5").Value = -0.0737094

$ws3.Range("A6").Value = "Z2"
$ws3.Range("B6").Value = 5.0370699999999999
$ws3.Range("C6").Value = 0.106265

$ws3.Range("A7").Value = "Z3"
$ws3.Range("B7").Value = 5.0641100000000003
$ws3.Range("C7").Value = 0.0017479100000000001

# The offsets are displayed with 12 decimal places
$ws3.Range("B2:C7").NumberFormat = "0.000000000000"

# Make Sheet3 the active sheet/tab and leave the selection on C7,
# matching the saved workbook view state.
$ws3.Activate()
$ws3.Range("C7").Select()
